$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.806.48"
$ws.Range("D3").Value = "2.303.20"
$ws.Range("E3").Value = "  -4.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.15"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.49"
$ws.Range("E6").Value = "  -3.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").Value = "2.302.95"
$ws.Range("E9").Value = "  -4.59%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  -5.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.02"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "2.719.03"
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("D16").Value = "58.775.04"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "2.302.02"
$ws.Range("E18").Value = "  -4.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -4.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("E20").Value = "  -4.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.16"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  -4.96%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.32"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.09"
$ws.Range("E27").Value = "  -6.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.32"
$ws.Range("E28").Value = "  -6.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.75"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.91"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "0.0₃0728"
$ws.Range("E31").Value = "  -5.77%  "
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("E33").Value = "  -5.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.383"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.76"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("E39").Value = "  -5.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.11"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -5.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "301.43"
$ws.Range("E42").Value = "  -7.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.40"
$ws.Range("E43").Value = "  -3.34%  "
$ws.Range("E44").Value = "  -5.41%  "
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0502"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.54"
$ws.Range("E48").Value = "  -6.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0215"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.63"
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("E51").Value = "  -0.33%  "
